# Generate Report for Archive
# Update localization status from "Ready for handoff" to "In Translation"
# across the Overview, zh-cn and de-de sheets, and re-pack the now-narrower
# Status columns to reflect the shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: status mirrored per-language in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$lastRow = $wsOverview.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    if ($wsOverview.Cells.Item($r, 5).Value2 -eq $oldStatus) {
        $wsOverview.Cells.Item($r, 5).Value = $newStatus
    }
    if ($wsOverview.Cells.Item($r, 6).Value2 -eq $oldStatus) {
        $wsOverview.Cells.Item($r, 6).Value = $newStatus
    }
}
# Columns narrow now that "In Translation" is shorter than "Ready for handoff"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-language detail sheets: status lives in column C ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 3).Value2 -eq $oldStatus) {
            $ws.Cells.Item($r, 3).Value = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
